$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "2016-03-21 10:40:21"
$wsZh.Range("H4").Value = "2016-03-21 10:40:47"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "2016-03-21 10:40:25"
$wsDe.Range("H4").Value = "2016-03-21 10:40:52"
